$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at EZ (column 156). This shifts the existing
# EZ ("nom") and FA ("url_produit") columns one place to the right,
# becoming FA and FB respectively, and leaves the new EZ column blank.
$ws.Range("EZ1").EntireColumn.Insert()

# New scrape timestamp header for the freshly inserted column.
$ws.Range("EZ1").Value = "2026-02-04 03:28:30"

# For the product rows that already had a price in column EY (the most
# recent prior scrape), carry that same price forward into the new EZ
# column (rows 2-80). Rows 81-206 had no price yet in EY, so their new
# EZ cell is left blank, matching the source data.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 155).Value2
    $ws.Cells.Item($r, 156).Value2 = $price
}
